$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param($ws, $startRow, $rows)
    for ($i = 0; $i -lt $rows.Count; $i++) {
        $row = $startRow + $i
        $vals = $rows[$i]
        for ($j = 0; $j -lt $vals.Count; $j++) {
            $cell = $ws.Cells.Item($row, $j + 1)
            $cell.NumberFormat = "@"
            $cell.Value = $vals[$j]
        }
    }
}

$pirRows = @(
    @("2026-01-28","16:32:46","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:32:49","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:32:51","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:32:55","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:33:00","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:33:05","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:33:10","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:33:15","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:33:20","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:33:25","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:33:30","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:33:35","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:33:40","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:33:46","16:00","Bathroom","No Motion","Inactive")
)

$humidityRows = @(
    @("2026-01-28","16:32:47","16:00","Bathroom","87.0%","Active"),
    @("2026-01-28","16:32:49","16:00","Bathroom","87.9%","Active"),
    @("2026-01-28","16:32:56","16:00","Bathroom","87.0%","Active"),
    @("2026-01-28","16:33:01","16:00","Bathroom","87.9%","Active"),
    @("2026-01-28","16:33:09","16:00","Bathroom","87.0%","Active"),
    @("2026-01-28","16:33:13","16:00","Bathroom","87.9%","Active"),
    @("2026-01-28","16:33:17","16:00","Bathroom","87.0%","Active"),
    @("2026-01-28","16:33:21","16:00","Bathroom","87.9%","Active"),
    @("2026-01-28","16:33:29","16:00","Bathroom","87.9%","Active"),
    @("2026-01-28","16:33:33","16:00","Bathroom","87.9%","Active"),
    @("2026-01-28","16:33:37","16:00","Bathroom","87.0%","Active"),
    @("2026-01-28","16:33:41","16:00","Bathroom","88.0%","Active"),
    @("2026-01-28","16:33:45","16:00","Bathroom","87.9%","Active")
)

$temperatureRows = @(
    @("2026-01-28","16:32:48","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:32:50","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:32:57","16:00","Bathroom","22.9C","Active"),
    @("2026-01-28","16:33:02","16:00","Bathroom","22.9C","Active"),
    @("2026-01-28","16:33:09","16:00","Bathroom","22.9C","Active"),
    @("2026-01-28","16:33:13","16:00","Bathroom","22.9C","Active"),
    @("2026-01-28","16:33:17","16:00","Bathroom","22.9C","Active"),
    @("2026-01-28","16:33:22","16:00","Bathroom","22.9C","Active"),
    @("2026-01-28","16:33:30","16:00","Bathroom","22.9C","Active"),
    @("2026-01-28","16:33:34","16:00","Bathroom","22.9C","Active"),
    @("2026-01-28","16:33:37","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:33:42","16:00","Bathroom","22.9C","Active"),
    @("2026-01-28","16:33:45","16:00","Bathroom","22.8C","Active")
)

$wsPIR = $wb.Worksheets.Item("PIR")
Add-LogRows $wsPIR 271 $pirRows

$wsHumidity = $wb.Worksheets.Item("Humidity")
Add-LogRows $wsHumidity 259 $humidityRows

$wsTemperature = $wb.Worksheets.Item("Temperature")
Add-LogRows $wsTemperature 260 $temperatureRows
